$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append -- fill column by column (year, publisher,
# image, japanese, english, product_type) to match the order in which the
# author actually typed the values (and thus the shared-string table order).
$ws.Range("A4").Value = 2006
$ws.Range("A5").Value = 2008
$ws.Range("A6").Value = 2009

$ws.Range("D4").Value = "SoftBank Creative"
$ws.Range("D5").Value = "SoftBank Creative"
$ws.Range("D6").Value = "SoftBank Creative"

$ws.Range("E4").Value = "sword_and_magic_museum.jpg"
$ws.Range("E5").Value = "monster_edition.jpg"
$ws.Range("E6").Value = "monster_edition_vol2.jpg"

$ws.Range("B4").Value = "新説 RPG幻想事典 剣と魔法の博物誌"
$ws.Range("B5").Value = "新説RPG幻想事典 剣と魔法の博物誌~モンスター編~"
$ws.Range("B6").Value = "新説 RPG幻想事典 剣と魔法の博物誌 ~モンスター編2~"

$ws.Range("C4").Value = "New RPG Fantasy Encyclopedia Sword and Magic Museum"
$ws.Range("C5").Value = "New RPG Fantasy Encyclopedia Sword and Magic Museum: Monster Edition"
$ws.Range("C6").Value = "New RPG Illusion Encyclopedia Sword and Magic Museum: Monster Edition 2"

$ws.Range("F4").Value = "supplement"
$ws.Range("F5").Value = "supplement"
$ws.Range("F6").Value = "supplement"

# Bold the header row
$headerRow = $ws.Range("A1:F1")
$ws.Rows.Item(1).Select() | Out-Null
$headerRow.Font.Bold = $true

# Widen columns B and C to fit the new, longer content
$ws.Columns.Item(2).ColumnWidth = 53.666666666666664
$ws.Columns.Item(3).ColumnWidth = 65.833333333333333
